$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.360.75"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.858.62"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'314.40"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.4638"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "'0.3719"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "'0.07346"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'0.8832"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "'0.07895"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").Value = "'19.89"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.891.53"
$ws.Range("E13").Value = "  +4.76%  "
$ws.Range("D14").Value = "'5.398"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "'6.577"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'92.09"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'0.000008869"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "27.400.27"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "'5.132"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.088.75"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'1.892"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'153.07"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "'2.084"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").Value = "'5.129"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'116.41"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'0.08892"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'3.027"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Value = "'0.7567"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("D34").Value = "'1.162"
$ws.Range("D35").Value = "'4.489"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").Value = "'2.619"
$ws.Range("E36").Value = "  +8.87%  "
$ws.Range("D37").Value = "'0.01960"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").Value = "'1.078"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05234"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.973"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "'7.110"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'0.5169"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'0.1645"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'8.341"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "'0.4835"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'10.30"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "'103.65"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'1.653"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'0.06238"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'65.71"
$ws.Range("E51").Value = "  +2.49%  "
